$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRange, $text) {
    $cellRange.NumberFormat = "@"
    $cellRange.Value = $text
    $cellRange.Style = "Normal"
}

# Row 30 / 31 swap: Filecoin <-> EthereumClassic
$ws.Range("B30").Value = "EthereumClassic"
$ws.Range("C30").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
Set-TextValue $ws.Range("D30") "30.78"
$ws.Range("E30").Value = "  +1.96%  "

$ws.Range("B31").Value = "Filecoin"
$ws.Range("C31").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
Set-TextValue $ws.Range("D31") "8.45"
$ws.Range("E31").Value = "  -2.51%  "

$ws.Range("D2").Value = "66.238.98"
$ws.Range("E2").Value = "  +1.30%  "
$ws.Range("D3").Value = "3.341.03"
$ws.Range("E3").Value = "  +2.77%  "
Set-TextValue $ws.Range("D4") "0.999"
$ws.Range("E4").Value = "  -0.27%  "
Set-TextValue $ws.Range("D5") "189.11"
$ws.Range("E5").Value = "  +4.84%  "
Set-TextValue $ws.Range("D6") "556.51"
$ws.Range("E6").Value = "  +0.20%  "
Set-TextValue $ws.Range("D7") "1.00"
$ws.Range("E7").Value = "  +0.04%  "
$ws.Range("D8").Value = "3.335.05"
$ws.Range("E8").Value = "  +2.83%  "
Set-TextValue $ws.Range("D9") "0.582"
$ws.Range("E9").Value = "  -1.46%  "
Set-TextValue $ws.Range("D10") "0.178"
Set-TextValue $ws.Range("D11") "0.581"
$ws.Range("E11").Value = "  -0.50%  "
Set-TextValue $ws.Range("D12") "46.28"
$ws.Range("E12").Value = "  -1.75%  "
Set-TextValue $ws.Range("D13") "0.0000268"
$ws.Range("E13").Value = "  +1.45%  "
$ws.Range("D14").Value = "3.873.93"
$ws.Range("E14").Value = "  +2.36%  "
Set-TextValue $ws.Range("D15") "8.52"
$ws.Range("E15").Value = "  -0.20%  "
Set-TextValue $ws.Range("D16") "586.38"
$ws.Range("E16").Value = "  -7.45%  "
$ws.Range("D17").Value = "66.229.11"
$ws.Range("E17").Value = "  +1.25%  "
$ws.Range("D18").Value = "3.348.97"
$ws.Range("E18").Value = "  +2.74%  "
$ws.Range("E19").Value = "  +1.14%  "
Set-TextValue $ws.Range("D20") "17.88"
$ws.Range("E20").Value = "  +1.45%  "
Set-TextValue $ws.Range("D21") "10.97"
$ws.Range("E21").Value = "  -3.03%  "
Set-TextValue $ws.Range("D22") "0.897"
$ws.Range("E22").Value = "  -0.10%  "
Set-TextValue $ws.Range("D23") "18.22"
$ws.Range("E23").Value = "  +2.74%  "
Set-TextValue $ws.Range("D24") "5.02"
$ws.Range("E24").Value = "  +1.91%  "
Set-TextValue $ws.Range("D25") "99.37"
$ws.Range("E25").Value = "  -5.14%  "
$ws.Range("E26").Value = "  +0.05%  "
$ws.Range("E27").Value = "  +0.74%  "
Set-TextValue $ws.Range("D28") "2.70"
$ws.Range("E28").Value = "  +2.05%  "
Set-TextValue $ws.Range("D29") "9.41"
$ws.Range("E29").Value = "  -0.53%  "
Set-TextValue $ws.Range("D32") "6.64"
$ws.Range("E32").Value = "  +5.99%  "
Set-TextValue $ws.Range("D33") "3.79"
$ws.Range("E33").Value = "  -5.39%  "
Set-TextValue $ws.Range("D34") "579.08"
$ws.Range("E34").Value = "  +6.20%  "
Set-TextValue $ws.Range("D35") "10.91"
$ws.Range("E35").Value = "  -0.71%  "
$ws.Range("E36").Value = "  -0.87%  "
$ws.Range("D37").Value = "3.761.69"
$ws.Range("E37").Value = "  +4.98%  "
$ws.Range("E38").Value = "  -0.08%  "
Set-TextValue $ws.Range("D39") "56.33"
$ws.Range("E39").Value = "  -0.94%  "
Set-TextValue $ws.Range("D40") "34.49"
$ws.Range("E40").Value = "  +8.35%  "
$ws.Range("D41").Value = "0.0₃0698"
$ws.Range("E41").Value = "  -2.33%  "
Set-TextValue $ws.Range("D42") "0.126"
$ws.Range("E42").Value = "  -1.63%  "
Set-TextValue $ws.Range("D43") "2.64"
$ws.Range("E43").Value = "  -2.54%  "
$ws.Range("E44").Value = "  -6.63%  "
Set-TextValue $ws.Range("D45") "3.37"
$ws.Range("E45").Value = "  +1.67%  "
Set-TextValue $ws.Range("D46") "0.337"
$ws.Range("E46").Value = "  +1.05%  "
Set-TextValue $ws.Range("D47") "0.0412"
$ws.Range("E47").Value = "  -0.63%  "
Set-TextValue $ws.Range("D48") "3.02"
$ws.Range("E48").Value = "  -10.50%  "
Set-TextValue $ws.Range("D49") "0.128"
$ws.Range("E49").Value = "  -0.01%  "
$ws.Range("E50").Value = "  +0.05%  "
Set-TextValue $ws.Range("D51") "2.54"
$ws.Range("E51").Value = "  -1.96%  "
